$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Euclides-Usinagem -> -
$ws.Range("D3").Value = "-"

# Row 4: Jorge Aquino-Motores de aplicação -> -
$ws.Range("C4").Value = "-"

# Row 6: - -> Euclides-Usinagem
$ws.Range("D6").Value = "Euclides-Usinagem"

# Row 7: - -> Jorge Aquino-Motores de aplicação
$ws.Range("C7").Value = "Jorge Aquino-Motores de aplicação"

# Row 8: Almoço -> - (all columns B:F)
$ws.Range("B8:F8").Value = "-"

# Row 9: 13:00 -> 12:20, and - -> Almoço (B:F)
$ws.Range("A9").Value = "12:20"
$ws.Range("B9:F9").Value = "Almoço"

# Row 10: 13:50 -> 13:00
$ws.Range("A10").Value = "13:00"

# Row 11: 14:40 -> 13:50
$ws.Range("A11").Value = "13:50"

# Row 12: 15:30 -> 14:40, Intervalo -> -
$ws.Range("A12").Value = "14:40"
$ws.Range("B12:F12").Value = "-"

# Row 13: 15:50 -> 15:30, - -> Intervalo
$ws.Range("A13").Value = "15:30"
$ws.Range("B13:F13").Value = "Intervalo"

# Row 14 (new): 15:50, all -
$ws.Range("A14").Value = "15:50"
$ws.Range("B14:F14").Value = "-"

# Row 15 (new, was old row 14 content): 16:40, all -
$ws.Range("A15").Value = "16:40"
$ws.Range("B15:F15").Value = "-"

# Row 16 (new): 17:30, all -
$ws.Range("A16").Value = "17:30"
$ws.Range("B16:F16").Value = "-"

# Row 17 (new): 18:20, B:F empty (touch formatting so the empty cells are materialized)
$ws.Range("A17").Value = "18:20"
$ws.Range("B17:F17").Font.Bold = $false
